$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C. This shifts the existing
# column C ("Jun_10" header + "UN"/special values) to column E.
$ws.Range("C1:D1").EntireColumn.Insert()

# Match the column width formatting used elsewhere on the sheet (8.0
# character width, custom width flag set) for the two newly inserted
# columns, and re-assert it on the shifted former column C (now E)
# since the column insert drops its custom-width flag.
$ws.Columns("C").ColumnWidth = 7.1875
$ws.Columns("D").ColumnWidth = 7.1875
$ws.Columns("E").ColumnWidth = 7.1875

# Row 1 headers: newest date goes in B, then each successive column
# holds an older date, oldest ("Jun_10") already sits in E after the
# insert/shift above.
$ws.Cells.Item(1, 2).Value = "Jun_17"
$ws.Cells.Item(1, 3).Value = "Jun_15"
$ws.Cells.Item(1, 4).Value = "Jun_13"

# Fill the two new data columns (C, D) with the same "UN" placeholder
# used throughout column B/E for every data row.
For ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}
